# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Updates the DAMSLTag (column I) and DialogAct (column J) values for a set of rows
# in the dialog-act annotated transcript sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of row number -> @(new DAMSLTag, new DialogAct)
$changes = @{
    30  = @("b",  "Acknowledge (Backchannel)")
    39  = @("b",  "Acknowledge (Backchannel)")
    53  = @("aa", "Agree/Accept")
    77  = @("sd", "Statement-non-opinion")
    79  = @("sd", "Statement-non-opinion")
    93  = @("ba", "Appreciation")
    115 = @("sd", "Statement-non-opinion")
    116 = @("sd", "Statement-non-opinion")
    125 = @("b",  "Acknowledge (Backchannel)")
    137 = @("aa", "Agree/Accept")
    148 = @("sd", "Statement-non-opinion")
    149 = @("sv", "Statement-opinion")
    164 = @("%",  "Uninterpretable")
    181 = @("sd", "Statement-non-opinion")
    189 = @("aa", "Agree/Accept")
    201 = @("aa", "Agree/Accept")
    211 = @("sd", "Statement-non-opinion")
    218 = @("sd", "Statement-non-opinion")
    219 = @("b",  "Acknowledge (Backchannel)")
    251 = @("ba", "Appreciation")
    261 = @("sv", "Statement-opinion")
    266 = @("sv", "Statement-opinion")
    269 = @("b",  "Acknowledge (Backchannel)")
    296 = @("sd", "Statement-non-opinion")
    297 = @("sv", "Statement-opinion")
    298 = @("sd", "Statement-non-opinion")
    303 = @("aa", "Agree/Accept")
    304 = @("aa", "Agree/Accept")
    307 = @("aa", "Agree/Accept")
    312 = @("sd", "Statement-non-opinion")
    326 = @("ba", "Appreciation")
    337 = @("sd", "Statement-non-opinion")
    345 = @("aa", "Agree/Accept")
    346 = @("aa", "Agree/Accept")
    348 = @("sv", "Statement-opinion")
    355 = @("aa", "Agree/Accept")
    366 = @("%",  "Uninterpretable")
    396 = @("aa", "Agree/Accept")
    397 = @("b",  "Acknowledge (Backchannel)")
    413 = @("sd", "Statement-non-opinion")
}

foreach ($rowNum in $changes.Keys) {
    $values = $changes[$rowNum]
    $ws.Cells.Item($rowNum, 9).Value = $values[0]
    $ws.Cells.Item($rowNum, 10).Value = $values[1]
}
